$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.367.69"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "1.941.40"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.20"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7321"
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3346"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07331"
$ws.Range("E10").Value = "  +4.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8166"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08107"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.937.27"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.493"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.11"
$ws.Range("E15").Value = "  -5.44%  "
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "30.368.86"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008314"
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.14"
$ws.Range("E19").Value = "  -7.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.922"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "2.192.86"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.975"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.837"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.69"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.407"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.43"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1330"
$ws.Range("E29").Value = "  -9.11%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.349"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.465"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.231"
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.312"
$ws.Range("E35").Value = "  +7.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7549"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.739"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01982"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.849"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.65"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.562"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.050"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8476"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.69"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.860"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.503"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.01"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4216"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06064"
$ws.Range("E51").Value = "  +1.45%  "
